$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old column D ("mx1" header + value) entirely, shifting cells left.
$ws.Range("D1:D2").Delete(-4121)

# Update header row: B1/C1 keep their bold/bordered style, just change the text.
$ws.Range("B1").Value = "Condition"
$ws.Range("C1").Value = "ddCT"

# Existing row 2 (style 1 on A2) becomes the "dusp11" row; keep A2 = 0.
$ws.Range("B2").Value = "dusp11"
$ws.Range("C2").Value = 0.5864994333333335

# Insert two new rows below row 2 for the "ifnb" / "mx1" entries, pushed from
# row 3 onward (sheet only has 2 rows so this just extends it).
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "ifnb"
$ws.Range("C3").Value = -1.147916233333333

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "mx1"
$ws.Range("C4").Value = -0.2273942333333316

# Give the new A3/A4 cells the same style (bold + border) as A2.
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
